$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Kanban Board sheet: clear its "tabSelected" / update selection to the
#    full data range (it will no longer be the active tab once we select
#    the Scrum Board sheet further down).
# ---------------------------------------------------------------------
$wsKanban = $wb.Worksheets.Item("Kanban Board")
$wsKanban.Range("A1:G13").Select()

# ---------------------------------------------------------------------
# 2) Scrum Board sheet: populate the Scrum test-plan rows (mirrors the
#    Kanban Board sheet layout/structure with Scrum-specific content).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Scrum Board")

# Header row
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Area"
$ws.Range("C1").Value = "Test case"
$ws.Range("D1").Value = "Pre-requisite"
$ws.Range("E1").Value = "Steps"
$ws.Range("F1").Value = "Outcome"
$ws.Range("G1").Value = "Comments"

# Row 2 - SB-1
$ws.Range("A2").Value = "SB-1"
$ws.Range("B2").Value = "Sidebar"
$ws.Range("C2").Value = "Initial screen"
$ws.Range("D2").Value = "Projects in WD starting with ""Scrum"",`nScrum config completed,`nLogin to WD tenant WCP-\enabled"
$ws.Range("F2").Value = "Defined scrum teams are listed in the sidebar"

# Row 3 - SB-2
$ws.Range("A3").Value = "SB-2"
$ws.Range("B3").Value = "Sidebar"
$ws.Range("C3").Value = "Refresh Team list"
$ws.Range("E3").Value = "Push ""refresh team list"" button"
$ws.Range("F3").Value = "Team list is loaded"

# Row 4 - SB-3
$ws.Range("A4").Value = "SB-3"
$ws.Range("B4").Value = "Sidebar"
$ws.Range("C4").Value = "Load team"
$ws.Range("E4").Value = "Select a team on the scrum list"
$ws.Range("F4").Value = "Team name is shown`nNew task and Refresh buttons are available`nProgress bar is shown`nSprint list is shown (default to current)`nTasks are loaded in Board`nColumn summaries are updated"

# Row 5 - SB-4
$ws.Range("A5").Value = "SB-4"
$ws.Range("B5").Value = "Board"
$ws.Range("C5").Value = "Drag task"
$ws.Range("D5").Value = "Board loaded with project tasks"

# Row 6 - SB-5
$ws.Range("A6").Value = "SB-5"
$ws.Range("B6").Value = "Board"
$ws.Range("C6").Value = "Edit task - update fields"
$ws.Range("D6").Value = "Board loaded with project tasks"

# Row 7 - SB-6
$ws.Range("A7").Value = "SB-6"
$ws.Range("B7").Value = "Board"
$ws.Range("C7").Value = "Edit task - add resource"
$ws.Range("F7").Value = "Available resources are shown,`nnew assignment is persisted"

# Row 8 - SB-7
$ws.Range("A8").Value = "SB-7"
$ws.Range("B8").Value = "Board"
$ws.Range("C8").Value = "Move sprint"
$ws.Range("F8").Value = "Task is moved to a new sprint"

# Row 9 - SB-8
$ws.Range("A9").Value = "SB-8"
$ws.Range("B9").Value = "Board"
$ws.Range("C9").Value = "New Task"
$ws.Range("D9").Value = "Board loaded with project tasks"

# Row 10 - SB-9
$ws.Range("A10").Value = "SB-9"
$ws.Range("B10").Value = "Board"
$ws.Range("C10").Value = "Refresh"
$ws.Range("D10").Value = "Board loaded with project tasks"

# Row 11 - SB-10
$ws.Range("A11").Value = "SB-10"
$ws.Range("B11").Value = "Navigation"
$ws.Range("C11").Value = "Navigate config"

# Row 12 - SB-11
$ws.Range("A12").Value = "SB-11"
$ws.Range("B12").Value = "Navigation"
$ws.Range("C12").Value = "Navigate back"

# Row 13 - SB-12
$ws.Range("A13").Value = "SB-12"
$ws.Range("B13").Value = "Navigation"
$ws.Range("C13").Value = "Search"

# Row 14 - SB-13
$ws.Range("A14").Value = "SB-13"
$ws.Range("B14").Value = "Navigation"
$ws.Range("C14").Value = "Expand/Hide sidebar"

# Formatting: vertical-top alignment across the whole table, matching the
# Kanban Board sheet's style, and wrap-text on the long multi-line cells.
$ws.Range("A1:G14").VerticalAlignment = -4160
$ws.Range("D2").WrapText = $true
$ws.Range("F4").WrapText = $true
$ws.Range("F7").WrapText = $true

# Row heights for the wrapped multi-line cells (matches Kanban Board sheet).
$ws.Rows(2).RowHeight = 43.5
$ws.Rows(4).RowHeight = 116
$ws.Rows(7).RowHeight = 29

# Column widths (A-C best-fit to content, D-F custom) matching target layout.
$ws.Columns(1).ColumnWidth = 4.8327
$ws.Columns(2).ColumnWidth = 9.333
$ws.Columns(3).ColumnWidth = 20.333
$ws.Columns(4).ColumnWidth = 31.6673
$ws.Columns(5).ColumnWidth = 29.666
$ws.Columns(6).ColumnWidth = 32.1673

# Activate the Scrum Board sheet and select C14 - it becomes the workbook's
# active tab/sheet (activeTab), with the Kanban Board sheet no longer
# selected.
$ws.Range("C14").Select()
